$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = 0.5
$ws.Range("C18").Value = "Implementation+ Testing"

$ws.Range("B20").Select()
